$d = $word.ActiveDocument

# Locate the paragraph "Look into how Process handle is managed..." — the
# last bullet in the "General" sub-list before the "New Modules" heading.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Look into how Process handle is managed*") {
        $target = $p
    }
}

# Insert a new paragraph after it and fill it in (inherits the ListParagraph
# style / numbering of $target automatically).
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p1 = $target.Next()
$p1.Range.Text = "Support compiling HadesMem as a DLL."

# Insert another new paragraph after that one.
$r2 = $p1.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.Text = "Change to use shared CRT linking."

# Word keeps the hidden "_GoBack" bookmark at the location of the most
# recent edit; move it here to match.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$endRange = $p2.Range
$endRange.Collapse(0)
$endRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $endRange)

$d.Save()
